# Applies the 2024-10-18 data update to the Chicago violent crime workbook.
# For each affected worksheet, updates the 2024 (column K) values for the
# specified rows to reflect newly reported incidents.

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Citywide Totals")
$ws.Cells.Item(2, 11).Value = 6179
$ws.Cells.Item(3, 11).Value = 6366
$ws.Cells.Item(4, 11).Value = 1334
$ws.Cells.Item(5, 11).Value = 452
$ws.Cells.Item(6, 11).Value = 7014
$ws.Cells.Item(7, 11).Value = 21345

$ws = $wb.Worksheets.Item("Logan Square")
$ws.Cells.Item(3, 11).Value = 72
$ws.Cells.Item(7, 11).Value = 271

$ws = $wb.Worksheets.Item("Austin")
$ws.Cells.Item(3, 11).Value = 427
$ws.Cells.Item(4, 11).Value = 76
$ws.Cells.Item(6, 11).Value = 472
$ws.Cells.Item(7, 11).Value = 1397

$ws = $wb.Worksheets.Item("South Chicago")
$ws.Cells.Item(2, 11).Value = 161
$ws.Cells.Item(7, 11).Value = 466

$ws = $wb.Worksheets.Item("Garfield Park")
$ws.Cells.Item(3, 11).Value = 334
$ws.Cells.Item(7, 11).Value = 929

$ws = $wb.Worksheets.Item("West Pullman")
$ws.Cells.Item(6, 11).Value = 80
$ws.Cells.Item(7, 11).Value = 352

$ws = $wb.Worksheets.Item("Grand Crossing")
$ws.Cells.Item(3, 11).Value = 242
$ws.Cells.Item(7, 11).Value = 722

$ws = $wb.Worksheets.Item("New City")
$ws.Cells.Item(2, 11).Value = 169
$ws.Cells.Item(7, 11).Value = 502

$ws = $wb.Worksheets.Item("Woodlawn")
$ws.Cells.Item(3, 11).Value = 144
$ws.Cells.Item(4, 11).Value = 17
$ws.Cells.Item(7, 11).Value = 351

$ws = $wb.Worksheets.Item("By Neighborhood")
$ws.Cells.Item(2, 11).Value = 188
$ws.Cells.Item(7, 11).Value = 625
$ws.Cells.Item(8, 11).Value = 1397
$ws.Cells.Item(10, 11).Value = 122
$ws.Cells.Item(11, 11).Value = 397
$ws.Cells.Item(15, 11).Value = 218
$ws.Cells.Item(19, 11).Value = 624
$ws.Cells.Item(20, 11).Value = 509
$ws.Cells.Item(21, 11).Value = 69
$ws.Cells.Item(23, 11).Value = 218
$ws.Cells.Item(24, 11).Value = 63
$ws.Cells.Item(25, 11).Value = 103
$ws.Cells.Item(27, 11).Value = 204
$ws.Cells.Item(29, 11).Value = 1162
$ws.Cells.Item(31, 11).Value = 237
$ws.Cells.Item(33, 11).Value = 929
$ws.Cells.Item(34, 11).Value = 122
$ws.Cells.Item(37, 11).Value = 722
$ws.Cells.Item(42, 11).Value = 792
$ws.Cells.Item(44, 11).Value = 180
$ws.Cells.Item(45, 11).Value = 29
$ws.Cells.Item(47, 11).Value = 148
$ws.Cells.Item(48, 11).Value = 268
$ws.Cells.Item(51, 11).Value = 275
$ws.Cells.Item(53, 11).Value = 271
$ws.Cells.Item(56, 11).Value = 23
$ws.Cells.Item(63, 11).Value = 60
$ws.Cells.Item(65, 11).Value = 502
$ws.Cells.Item(66, 11).Value = 67
$ws.Cells.Item(67, 11).Value = 840
$ws.Cells.Item(71, 11).Value = 63
$ws.Cells.Item(72, 11).Value = 108
$ws.Cells.Item(73, 11).Value = 189
$ws.Cells.Item(78, 11).Value = 239
$ws.Cells.Item(79, 11).Value = 538
$ws.Cells.Item(83, 11).Value = 466
$ws.Cells.Item(84, 11).Value = 169
$ws.Cells.Item(85, 11).Value = 994
$ws.Cells.Item(88, 11).Value = 227
$ws.Cells.Item(90, 11).Value = 196
$ws.Cells.Item(91, 11).Value = 244
$ws.Cells.Item(95, 11).Value = 352
$ws.Cells.Item(96, 11).Value = 222
$ws.Cells.Item(97, 11).Value = 168
$ws.Cells.Item(99, 11).Value = 351
$ws.Cells.Item(101, 11).Value = 21345

$ws = $wb.Worksheets.Item("Gage Park")
$ws.Cells.Item(6, 11).Value = 81
$ws.Cells.Item(7, 11).Value = 237

$ws = $wb.Worksheets.Item("North Lawndale")
$ws.Cells.Item(2, 11).Value = 232
$ws.Cells.Item(3, 11).Value = 304
$ws.Cells.Item(7, 11).Value = 840

$ws = $wb.Worksheets.Item("South Deering")
$ws.Cells.Item(2, 11).Value = 58
$ws.Cells.Item(7, 11).Value = 169

$ws = $wb.Worksheets.Item("Englewood")
$ws.Cells.Item(4, 11).Value = 57
$ws.Cells.Item(6, 11).Value = 329
$ws.Cells.Item(7, 11).Value = 1162

$ws = $wb.Worksheets.Item("Lake View")
$ws.Cells.Item(5, 11).Value = 1
$ws.Cells.Item(7, 11).Value = 268

$ws = $wb.Worksheets.Item("Chatham")
$ws.Cells.Item(2, 11).Value = 186
$ws.Cells.Item(3, 11).Value = 189
$ws.Cells.Item(7, 11).Value = 624

$ws = $wb.Worksheets.Item("Irving Park")
$ws.Cells.Item(3, 11).Value = 47
$ws.Cells.Item(7, 11).Value = 180

$ws = $wb.Worksheets.Item("Humboldt Park")
$ws.Cells.Item(3, 11).Value = 238
$ws.Cells.Item(6, 11).Value = 296
$ws.Cells.Item(7, 11).Value = 792

$ws = $wb.Worksheets.Item("Avondale")
$ws.Cells.Item(6, 11).Value = 56
$ws.Cells.Item(7, 11).Value = 122

$ws = $wb.Worksheets.Item("Rogers Park")
$ws.Cells.Item(2, 11).Value = 71
$ws.Cells.Item(7, 11).Value = 239

$ws = $wb.Worksheets.Item("Dunning")
$ws.Cells.Item(2, 11).Value = 25
$ws.Cells.Item(7, 11).Value = 63

$ws = $wb.Worksheets.Item("Douglas")
$ws.Cells.Item(3, 11).Value = 77
$ws.Cells.Item(7, 11).Value = 218

$ws = $wb.Worksheets.Item("West Ridge")
$ws.Cells.Item(2, 11).Value = 70
$ws.Cells.Item(7, 11).Value = 222

$ws = $wb.Worksheets.Item("Washington Park")
$ws.Cells.Item(6, 11).Value = 51
$ws.Cells.Item(7, 11).Value = 244

$ws = $wb.Worksheets.Item("Chinatown")
$ws.Cells.Item(6, 11).Value = 41
$ws.Cells.Item(7, 11).Value = 69

$ws = $wb.Worksheets.Item("Roseland")
$ws.Cells.Item(2, 11).Value = 179
$ws.Cells.Item(4, 11).Value = 33
$ws.Cells.Item(5, 11).Value = 17
$ws.Cells.Item(7, 11).Value = 538

$ws = $wb.Worksheets.Item("Chicago Lawn")
$ws.Cells.Item(3, 11).Value = 165
$ws.Cells.Item(7, 11).Value = 509

$ws = $wb.Worksheets.Item("Auburn Gresham")
$ws.Cells.Item(2, 11).Value = 204
$ws.Cells.Item(3, 11).Value = 207
$ws.Cells.Item(6, 11).Value = 169
$ws.Cells.Item(7, 11).Value = 625

$ws = $wb.Worksheets.Item("Garfield Ridge")
$ws.Cells.Item(3, 11).Value = 33
$ws.Cells.Item(7, 11).Value = 122

$ws = $wb.Worksheets.Item("East Side")
$ws.Cells.Item(3, 11).Value = 36
$ws.Cells.Item(7, 11).Value = 103

$ws = $wb.Worksheets.Item("Kenwood")
$ws.Cells.Item(6, 11).Value = 49
$ws.Cells.Item(7, 11).Value = 148

$ws = $wb.Worksheets.Item("Brighton Park")
$ws.Cells.Item(3, 11).Value = 55
$ws.Cells.Item(7, 11).Value = 218

$ws = $wb.Worksheets.Item("North Center")
$ws.Cells.Item(3, 11).Value = 17
$ws.Cells.Item(7, 11).Value = 67

$ws = $wb.Worksheets.Item("Belmont Cragin")
$ws.Cells.Item(6, 11).Value = 128
$ws.Cells.Item(7, 11).Value = 397

$ws = $wb.Worksheets.Item("Portage Park")
$ws.Cells.Item(3, 11).Value = 51
$ws.Cells.Item(7, 11).Value = 189

$ws = $wb.Worksheets.Item("Albany Park")
$ws.Cells.Item(3, 11).Value = 50
$ws.Cells.Item(7, 11).Value = 188

$ws = $wb.Worksheets.Item("West Town")
$ws.Cells.Item(6, 11).Value = 94
$ws.Cells.Item(7, 11).Value = 168

$ws = $wb.Worksheets.Item("United Center")
$ws.Cells.Item(6, 11).Value = 95
$ws.Cells.Item(7, 11).Value = 227

$ws = $wb.Worksheets.Item("Edgewater")
$ws.Cells.Item(2, 11).Value = 55
$ws.Cells.Item(6, 11).Value = 72
$ws.Cells.Item(7, 11).Value = 204

$ws = $wb.Worksheets.Item("Washington Heights")
$ws.Cells.Item(3, 11).Value = 56
$ws.Cells.Item(7, 11).Value = 196

$ws = $wb.Worksheets.Item("Little Italy, UIC")
$ws.Cells.Item(3, 11).Value = 74
$ws.Cells.Item(6, 11).Value = 89
$ws.Cells.Item(7, 11).Value = 275

$ws = $wb.Worksheets.Item("South Shore")
$ws.Cells.Item(2, 11).Value = 324
$ws.Cells.Item(3, 11).Value = 342
$ws.Cells.Item(5, 11).Value = 29
$ws.Cells.Item(6, 11).Value = 244
$ws.Cells.Item(7, 11).Value = 994

$ws = $wb.Worksheets.Item("Oakland")
$ws.Cells.Item(6, 11).Value = 17
$ws.Cells.Item(7, 11).Value = 63

$ws = $wb.Worksheets.Item("Old Town")
$ws.Cells.Item(2, 11).Value = 22
$ws.Cells.Item(3, 11).Value = 26
$ws.Cells.Item(7, 11).Value = 108

$ws = $wb.Worksheets.Item("Jackson Park")
$ws.Cells.Item(2, 11).Value = 7
$ws.Cells.Item(7, 11).Value = 29

$ws = $wb.Worksheets.Item("Magnificent Mile")
$ws.Cells.Item(6, 11).Value = 11
$ws.Cells.Item(7, 11).Value = 23

Write-Host "Applied 2024-10-18 crime data update."
